$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 127; this shifts rows 127-136 down to 128-137,
# preserving all of their existing data.
$ws.Rows("127:127").Insert()

# Populate the newly inserted row 127 with the new weekly price entry.
$ws.Cells.Item(127, 1).Value = 4
$ws.Cells.Item(127, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(127, 3).Value = "Los Lagos"
$ws.Cells.Item(127, 4).Value = 45013
$ws.Cells.Item(127, 5).Value = 10
$ws.Cells.Item(127, 6).Value = 100112031
$ws.Cells.Item(127, 7).Value = "Poroto verde"
$ws.Cells.Item(127, 8).Value = "Magnum"
$ws.Cells.Item(127, 9).Value = "Primera"
$ws.Cells.Item(127, 10).Value = 40
$ws.Cells.Item(127, 11).Value = 30000
$ws.Cells.Item(127, 12).Value = 30000
$ws.Cells.Item(127, 13).Value = 30000
$ws.Cells.Item(127, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(127, 15).Value = "Región Metropolitana"
$ws.Cells.Item(127, 16).Value = 1200
$ws.Cells.Item(127, 17).Value = 25
$ws.Cells.Item(127, 18).Value = "Hortaliza"
